# Update "想去人数" (interested-count, column F) values across the four
# sheets of the 北京-漫展信息 workbook, per the upstream data refresh
# commit "Update gh-pages to output generated at 456a3b4".
#
# The same underlying events appear on multiple tabs (展览 / 演出 /
# 本地生活 feed into the combined 全部类型 tab), so each event's count is
# bumped in every sheet where it occurs.

$wb = $excel.ActiveWorkbook

# --- Worksheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 9915
$ws.Cells.Item(8, 6).Value = 397
$ws.Cells.Item(9, 6).Value = 444
$ws.Cells.Item(11, 6).Value = 221
$ws.Cells.Item(13, 6).Value = 499
$ws.Cells.Item(14, 6).Value = 12546
$ws.Cells.Item(20, 6).Value = 43
$ws.Cells.Item(24, 6).Value = 2748
$ws.Cells.Item(29, 6).Value = 2165
$ws.Cells.Item(30, 6).Value = 1069
$ws.Cells.Item(31, 6).Value = 4248
$ws.Cells.Item(32, 6).Value = 3764
$ws.Cells.Item(33, 6).Value = 780
$ws.Cells.Item(34, 6).Value = 2645
$ws.Cells.Item(35, 6).Value = 3078
$ws.Cells.Item(36, 6).Value = 57
$ws.Cells.Item(37, 6).Value = 1360
$ws.Cells.Item(39, 6).Value = 783
$ws.Cells.Item(40, 6).Value = 41
$ws.Cells.Item(41, 6).Value = 128
$ws.Cells.Item(43, 6).Value = 623
$ws.Cells.Item(45, 6).Value = 152

# --- Worksheet "演出" (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(5, 6).Value = 50
$ws.Cells.Item(11, 6).Value = 32
$ws.Cells.Item(13, 6).Value = 52
$ws.Cells.Item(15, 6).Value = 12
$ws.Cells.Item(18, 6).Value = 11

# --- Worksheet "本地生活" (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 61

# --- Worksheet "全部类型" (All types, combined feed) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(6, 6).Value = 9915
$ws.Cells.Item(8, 6).Value = 50
$ws.Cells.Item(11, 6).Value = 397
$ws.Cells.Item(12, 6).Value = 444
$ws.Cells.Item(14, 6).Value = 221
$ws.Cells.Item(15, 6).Value = 499
$ws.Cells.Item(16, 6).Value = 12546
$ws.Cells.Item(18, 6).Value = 61
$ws.Cells.Item(24, 6).Value = 2748
$ws.Cells.Item(28, 6).Value = 2165
$ws.Cells.Item(29, 6).Value = 1069
$ws.Cells.Item(30, 6).Value = 4248
$ws.Cells.Item(31, 6).Value = 3764
$ws.Cells.Item(32, 6).Value = 780
$ws.Cells.Item(33, 6).Value = 2645
$ws.Cells.Item(34, 6).Value = 3078
$ws.Cells.Item(35, 6).Value = 57
$ws.Cells.Item(36, 6).Value = 1360
$ws.Cells.Item(38, 6).Value = 783
$ws.Cells.Item(39, 6).Value = 41
$ws.Cells.Item(40, 6).Value = 128
$ws.Cells.Item(43, 6).Value = 623
$ws.Cells.Item(45, 6).Value = 152
